$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 45)
$ws.Range("D2").Value = (Get-Date -Year 2021 -Month 2 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M2").Value = 350
$ws.Range("N2").Value = 3500
$ws.Range("O2").Value = 4000
$ws.Range("P2").Value = 3750
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 1875

# Row 3 (was row 7)
$ws.Range("D3").Value = (Get-Date -Year 2020 -Month 12 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M3").Value = 300
$ws.Range("N3").Value = 3000
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 3000
$ws.Range("S3").Value = 1500

# Row 4 (was row 8)
$ws.Range("D4").Value = (Get-Date -Year 2020 -Month 12 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N4").Value = 3000
$ws.Range("O4").Value = 3000
$ws.Range("P4").Value = 3000
$ws.Range("S4").Value = 1500

# Row 5 (was row 2)
$ws.Range("D5").Value = (Get-Date -Year 2020 -Month 12 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M5").Value = 140
$ws.Range("O5").Value = 4500
$ws.Range("P5").Value = 4250
$ws.Range("S5").Value = 2125

# Row 6 (was row 37)
$ws.Range("D6").Value = (Get-Date -Year 2021 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M6").Value = 150
$ws.Range("R6").Value = 'Provincia de Curicó'

# Row 7 (was row 38)
$ws.Range("D7").Value = (Get-Date -Year 2021 -Month 1 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 4000
$ws.Range("O7").Value = 4000
$ws.Range("P7").Value = 4000
$ws.Range("R7").Value = 'Provincia de Linares'
$ws.Range("S7").Value = 2000

# Row 8 (was row 35)
$ws.Range("D8").Value = (Get-Date -Year 2021 -Month 2 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M8").Value = 125
$ws.Range("N8").Value = 4000
$ws.Range("O8").Value = 4000
$ws.Range("P8").Value = 4000
$ws.Range("R8").Value = 'Provincia de Curicó'
$ws.Range("S8").Value = 2000

# Row 9 (was row 44)
$ws.Range("D9").Value = (Get-Date -Year 2020 -Month 12 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M9").Value = 100

# Row 10 (was row 20)
$ws.Range("D10").Value = (Get-Date -Year 2021 -Month 2 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M10").Value = 25

# Row 11 (was row 5)
$ws.Range("D11").Value = (Get-Date -Year 2020 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M11").Value = 250

# Row 12 (was row 34)
$ws.Range("D12").Value = (Get-Date -Year 2020 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M12").Value = 150
$ws.Range("R12").Value = 'Provincia de Linares'

# Row 13 (was row 27)
$ws.Range("D13").Value = (Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M13").Value = 250
$ws.Range("R13").Value = 'Provincia de Curicó'

# Row 14 (was row 28)
$ws.Range("D14").Value = (Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N14").Value = 4000
$ws.Range("O14").Value = 4000
$ws.Range("P14").Value = 4000
$ws.Range("S14").Value = 2000

# Row 15 (was row 10)
$ws.Range("D15").Value = (Get-Date -Year 2021 -Month 1 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 4000
$ws.Range("O15").Value = 4000
$ws.Range("P15").Value = 4000
$ws.Range("R15").Value = 'Provincia de Curicó'
$ws.Range("S15").Value = 2000

# Row 16 (was row 11)
$ws.Range("D16").Value = (Get-Date -Year 2021 -Month 1 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("R16").Value = 'Provincia de Linares'

# Row 17 (was row 30)
$ws.Range("D17").Value = (Get-Date -Year 2021 -Month 3 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M17").Value = 75
$ws.Range("R17").Value = 'Provincia de Curicó'

# Row 18 (was row 46)
$ws.Range("D18").Value = (Get-Date -Year 2021 -Month 1 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N18").Value = 3000
$ws.Range("O18").Value = 3500
$ws.Range("P18").Value = 3250
$ws.Range("S18").Value = 1625

# Row 19 (was row 26)
$ws.Range("D19").Value = (Get-Date -Year 2021 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N19").Value = 3000
$ws.Range("P19").Value = 3500
$ws.Range("S19").Value = 1750

# Row 20 (was row 43)
$ws.Range("D20").Value = (Get-Date -Year 2021 -Month 3 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M20").Value = 100

# Row 21 (was row 22)
$ws.Range("D21").Value = (Get-Date -Year 2021 -Month 1 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M21").Value = 150

# Row 22 (was row 23)
$ws.Range("M22").Value = 200
$ws.Range("R22").Value = 'Provincia de Linares'

# Row 23 (was row 41)
$ws.Range("D23").Value = (Get-Date -Year 2020 -Month 12 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M23").Value = 250

# Row 24 (was row 16)
$ws.Range("D24").Value = (Get-Date -Year 2020 -Month 12 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M24").Value = 300
$ws.Range("R24").Value = 'Provincia de Curicó'

# Row 25 (was row 17)
$ws.Range("D25").Value = (Get-Date -Year 2020 -Month 12 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M25").Value = 500
$ws.Range("N25").Value = 4000
$ws.Range("O25").Value = 4000
$ws.Range("P25").Value = 4000
$ws.Range("S25").Value = 2000

# Row 26 (was row 18)
$ws.Range("D26").Value = (Get-Date -Year 2021 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 4000
$ws.Range("P26").Value = 4000
$ws.Range("R26").Value = 'Provincia de Curicó'
$ws.Range("S26").Value = 2000

# Row 27 (was row 19)
$ws.Range("D27").Value = (Get-Date -Year 2021 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M27").Value = 400
$ws.Range("R27").Value = 'Provincia de Linares'

# Row 28 (was row 12)
$ws.Range("D28").Value = (Get-Date -Year 2021 -Month 1 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M28").Value = 250
$ws.Range("R28").Value = 'Provincia de Curicó'

# Row 29 (was row 13)
$ws.Range("D29").Value = (Get-Date -Year 2021 -Month 1 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M29").Value = 300

# Row 30 (was row 14)
$ws.Range("D30").Value = (Get-Date -Year 2020 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 3000
$ws.Range("O30").Value = 3000
$ws.Range("P30").Value = 3000
$ws.Range("R30").Value = 'Provincia de Linares'
$ws.Range("S30").Value = 1500

# Row 31 (was row 39)
$ws.Range("D31").Value = (Get-Date -Year 2021 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M31").Value = 170
$ws.Range("N31").Value = 3000
$ws.Range("O31").Value = 4000
$ws.Range("P31").Value = 3500
$ws.Range("R31").Value = 'Provincia de Linares'
$ws.Range("S31").Value = 1750

# Row 32 (was row 15)
$ws.Range("D32").Value = (Get-Date -Year 2020 -Month 12 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("N32").Value = 5000
$ws.Range("O32").Value = 5000
$ws.Range("P32").Value = 5000
$ws.Range("S32").Value = 2500

# Row 33 (was row 40)
$ws.Range("D33").Value = (Get-Date -Year 2021 -Month 1 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 750
$ws.Range("N33").Value = 4000
$ws.Range("O33").Value = 4000
$ws.Range("P33").Value = 4000
$ws.Range("R33").Value = 'Provincia de Curicó'
$ws.Range("S33").Value = 2000

# Row 34 (was row 42)
$ws.Range("D34").Value = (Get-Date -Year 2021 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M34").Value = 350

# Row 35 (was row 36)
$ws.Range("D35").Value = (Get-Date -Year 2021 -Month 2 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M35").Value = 75

# Row 36 (was row 3)
$ws.Range("D36").Value = (Get-Date -Year 2021 -Month 1 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M36").Value = 150

# Row 37 (was row 4)
$ws.Range("D37").Value = (Get-Date -Year 2021 -Month 1 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M37").Value = 250
$ws.Range("R37").Value = 'Provincia de Linares'

# Row 38 (was row 6)
$ws.Range("D38").Value = (Get-Date -Year 2020 -Month 12 -Day 21 -Hour 0 -Minute 0 -Second 0)

# Row 39 (was row 9)
$ws.Range("D39").Value = (Get-Date -Year 2021 -Month 1 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M39").Value = 200
$ws.Range("N39").Value = 4000
$ws.Range("P39").Value = 4000
$ws.Range("S39").Value = 2000

# Row 40 (was row 24)
$ws.Range("D40").Value = (Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M40").Value = 200
$ws.Range("R40").Value = 'Provincia de Linares'

# Row 41 (was row 31)
$ws.Range("D41").Value = (Get-Date -Year 2020 -Month 12 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M41").Value = 100
$ws.Range("N41").Value = 3400
$ws.Range("O41").Value = 3400
$ws.Range("P41").Value = 3400
$ws.Range("R41").Value = 'Provincia de Curicó'
$ws.Range("S41").Value = 1700

# Row 42 (was row 32)
$ws.Range("D42").Value = (Get-Date -Year 2020 -Month 12 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M42").Value = 200

# Row 43 (was row 33)
$ws.Range("D43").Value = (Get-Date -Year 2020 -Month 12 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L43").Value = 'Segunda'
$ws.Range("M43").Value = 50
$ws.Range("N43").Value = 3000
$ws.Range("O43").Value = 3000
$ws.Range("P43").Value = 3000
$ws.Range("R43").Value = 'Provincia de Linares'
$ws.Range("S43").Value = 1500

# Row 44 (was row 21)
$ws.Range("D44").Value = (Get-Date -Year 2021 -Month 2 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("R44").Value = 'Provincia de Curicó'

# Row 45 (was row 29)
$ws.Range("D45").Value = (Get-Date -Year 2021 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("M45").Value = 200
$ws.Range("N45").Value = 4000
$ws.Range("P45").Value = 4000
$ws.Range("R45").Value = 'Provincia de Linares'
$ws.Range("S45").Value = 2000

# Row 46 (was row 25)
$ws.Range("D46").Value = (Get-Date -Year 2020 -Month 12 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("O46").Value = 3000
$ws.Range("P46").Value = 3000
$ws.Range("R46").Value = 'Provincia de Linares'
$ws.Range("S46").Value = 1500

Write-Output "Applied reorder/update to 45 data rows"
